$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.697551846504211
$ws.Range("B1").Value = 1.954059481620789
$ws.Range("C1").Value = 2.008124828338623
$ws.Range("D1").Value = 2.572867870330811
$ws.Range("E1").Value = 3.454827070236206
